$d = $word.ActiveDocument

# -------------------------------------------------------------------
# 1) Merge the split/spell-checked runs for the title lines into
#    single runs (removes the w:proofErr wraps + joins the w:t's).
# -------------------------------------------------------------------
$d.Content.Find.Execute("Lappeenrannan teknillinen yliopisto", $true, $true, $false, $false, $false,
                         $true, 1, $false, "Lappeenrannan teknillinen yliopisto", 2)

$d.Content.Find.Execute("Sofware Development Skills", $true, $true, $false, $false, $false,
                         $true, 1, $false, "Sofware Development Skills", 2)

# -------------------------------------------------------------------
# 2) Merge the runs split around "a" / "url" (proofErr wraps) back
#    into a single run of plain text.
# -------------------------------------------------------------------
$d.Content.Find.Execute("Then deleted it and made a own as the code", $true, $true, $false, $false, $false,
                         $true, 1, $false, "Then deleted it and made a own as the code", 2)

$d.Content.Find.Execute("has the Id of map in url to look that specific map", $true, $true, $false, $false, $false,
                         $true, 1, $false, "has the Id of map in url to look that specific map", 2)

# -------------------------------------------------------------------
# 3) Append the new diary entries (29.3.2025 / 30.3.2025) at the end
#    of the document, after the trailing empty BodyText paragraph.
# -------------------------------------------------------------------

# Trailing empty paragraph already in the doc -> add "29.3.2025" para after it.
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "29.3.2025"

# "Some code refining..." paragraph.
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "Some code refining but didn’t succeed in making proper zoom"

# Blank separator paragraph - leave empty (no run).
$d.Paragraphs.Last.Range.InsertParagraphAfter()

# "30.3.2025" paragraph.
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "30.3.2025"

# Final paragraph with the drag-prevention note.
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "Prevented the dragging of an image causing download or opening a new tab when dragged outside of the website."
